$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove league weighting: set all competency weights (B2:B41) to 1
$ws.Range("B2:B41").Value = 1

# Update the view state to match the saved sheet view (scrolled/selected range)
$ws.Range("B2:B41").Select()
$excel.ActiveWindow.ScrollRow = 25
